$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 526, shifting existing rows 526:609 down to 527:610.
$ws.Rows.Item(526).Insert()

# Populate the newly inserted row with the new data point.
$ws.Range("A526").Value = 3
$ws.Range("B526").Value = "Femacal de La Calera"
$ws.Range("C526").Value = "Coquimbo"
$ws.Range("D526").Value = 45180
$ws.Range("E526").Value = 5
$ws.Range("F526").Value = 100112012
$ws.Range("G526").Value = "Espinaca"
$ws.Range("H526").Value = "Sin especificar"
$ws.Range("I526").Value = "Primera"
$ws.Range("J526").Value = 160
$ws.Range("K526").Value = 4000
$ws.Range("L526").Value = 4500
$ws.Range("M526").Value = 4188
$ws.Range("N526").Value = "$/docena de atados (3 kilos)"
$ws.Range("O526").Value = "Provincia de Quillota"
$ws.Range("P526").Value = 1396
$ws.Range("Q526").Value = 3
$ws.Range("R526").Value = "Hortaliza"
